{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same content changes described by the OOXML diff:\n//   1) Row 10 (\"\u0412\u0441\u0442\u0430\u0432\u0438\u0442\u044c \u043f\u043e\u0441\u043b\u0435 \u043a\u0430\u0436\u0434\u043e\u0433\u043e \u0441\u0442\u043e\u043b\u0431\u0446\u0430, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0435\u0433\u043e \u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439\n//      \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u0441\u0442\u0440\u043e\u043a\u0443 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\") gets reworded to\n//      \"...\u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439 \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0441\u0442\u043e\u043b\u0431\u0435\u0446 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\"\n//   2) Row 11 (\"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043d\u0435\u0433\u043e \u043a\u0430\u0436\u0434\u044b\u0439 \u0441\u0442\u043e\u043b\u0431\u0435\u0446, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0438\u0439 \u044d\u043b\u0435\u043c\u0435\u043d\u0442,\n//      \u043a\u0440\u0430\u0442\u043d\u044b\u0439 \u043f\u044f\u0442\u0438.\") gets reworded to \"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043c\u0430\u0441\u0441\u0438\u0432\u0430 \u043a\u0430\u0436\u0434\u044b\u0439\n//      \u0441\u0442\u043e\u043b\u0431\u0435\u0446...\" and the document's \"_GoBack\" bookmark moves to sit\n//      right after the newly typed word \"\u043c\u0430\u0441\u0441\u0438\u0432\u0430\" (simulating the\n//      cursor position Word leaves behind after the last edit).\n\nconst body = context.document.body;\n\n// --- Change 1: row 10 sentence rewording -------------------------------\nconst oldText10 =\n  \"\u0412\u0441\u0442\u0430\u0432\u0438\u0442\u044c \u043f\u043e\u0441\u043b\u0435 \u043a\u0430\u0436\u0434\u043e\u0433\u043e \u0441\u0442\u043e\u043b\u0431\u0446\u0430, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0435\u0433\u043e \u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439 \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u0441\u0442\u0440\u043e\u043a\u0443 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\";\nconst newText10 =\n  \"\u0412\u0441\u0442\u0430\u0432\u0438\u0442\u044c \u043f\u043e\u0441\u043b\u0435 \u043a\u0430\u0436\u0434\u043e\u0433\u043e \u0441\u0442\u043e\u043b\u0431\u0446\u0430, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0435\u0433\u043e \u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439 \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0441\u0442\u043e\u043b\u0431\u0435\u0446 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\";\n\nconst results10 = body.search(oldText10, { matchCase: true });\nresults10.load(\"items\");\nawait context.sync();\n\nif (results10.items.length > 0) {\n  results10.items[0].insertText(newText10, \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2: remove the pre-existing \"_GoBack\" bookmark --------------\n// Word keeps a single \"_GoBack\" bookmark that marks the location of the\n// last edit; it gets removed from its old spot and re-created at the new\n// edit location (see below).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Change 3: row 11 sentence rewording + relocate \"_GoBack\" ----------\nconst oldText11 =\n  \"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043d\u0435\u0433\u043e \u043a\u0430\u0436\u0434\u044b\u0439 \u0441\u0442\u043e\u043b\u0431\u0435\u0446, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0438\u0439 \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u043a\u0440\u0430\u0442\u043d\u044b\u0439 \u043f\u044f\u0442\u0438\";\nconst newText11 =\n  \"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043c\u0430\u0441\u0441\u0438\u0432\u0430 \u043a\u0430\u0436\u0434\u044b\u0439 \u0441\u0442\u043e\u043b\u0431\u0435\u0446, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0438\u0439 \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u043a\u0440\u0430\u0442\u043d\u044b\u0439 \u043f\u044f\u0442\u0438\";\n\nconst results11 = body.search(oldText11, { matchCase: true });\nresults11.load(\"items\");\nawait context.sync();\n\nif (results11.items.length > 0) {\n  results11.items[0].insertText(newText11, \"Replace\");\n  await context.sync();\n\n  // Find the freshly inserted word \"\u043c\u0430\u0441\u0441\u0438\u0432\u0430\" within this same sentence so\n  // we can drop the \"_GoBack\" bookmark immediately after it.\n  const afterReplace = body.search(\"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043c\u0430\u0441\u0441\u0438\u0432\u0430\", { matchCase: true });\n  afterReplace.load(\"items\");\n  await context.sync();\n\n  if (afterReplace.items.length > 0) {\n    const endRange = afterReplace.items[0].getRange(\"End\");\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same content changes described by the OOXML diff:\n#   1) Row 10 (\"\u0412\u0441\u0442\u0430\u0432\u0438\u0442\u044c \u043f\u043e\u0441\u043b\u0435 \u043a\u0430\u0436\u0434\u043e\u0433\u043e \u0441\u0442\u043e\u043b\u0431\u0446\u0430, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0435\u0433\u043e \u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439\n#      \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u0441\u0442\u0440\u043e\u043a\u0443 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\") gets reworded to\n#      \"...\u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439 \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0441\u0442\u043e\u043b\u0431\u0435\u0446 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\"\n#   2) Row 11 (\"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043d\u0435\u0433\u043e \u043a\u0430\u0436\u0434\u044b\u0439 \u0441\u0442\u043e\u043b\u0431\u0435\u0446, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0438\u0439 \u044d\u043b\u0435\u043c\u0435\u043d\u0442,\n#      \u043a\u0440\u0430\u0442\u043d\u044b\u0439 \u043f\u044f\u0442\u0438.\") gets reworded to \"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043c\u0430\u0441\u0441\u0438\u0432\u0430 \u043a\u0430\u0436\u0434\u044b\u0439\n#      \u0441\u0442\u043e\u043b\u0431\u0435\u0446...\" and the document's \"_GoBack\" bookmark moves to sit\n#      right after the newly typed word \"\u043c\u0430\u0441\u0441\u0438\u0432\u0430\" (simulating the\n#      cursor position Word leaves behind after the last edit).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: row 10 sentence rewording --------------------------------\n$range1 = $d.Content\n$found1 = $range1.Find.Execute(\n    \"\u0412\u0441\u0442\u0430\u0432\u0438\u0442\u044c \u043f\u043e\u0441\u043b\u0435 \u043a\u0430\u0436\u0434\u043e\u0433\u043e \u0441\u0442\u043e\u043b\u0431\u0446\u0430, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0435\u0433\u043e \u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439 \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u0441\u0442\u0440\u043e\u043a\u0443 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u0412\u0441\u0442\u0430\u0432\u0438\u0442\u044c \u043f\u043e\u0441\u043b\u0435 \u043a\u0430\u0436\u0434\u043e\u0433\u043e \u0441\u0442\u043e\u043b\u0431\u0446\u0430, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0435\u0433\u043e \u043c\u0430\u043a\u0441\u0438\u043c\u0430\u043b\u044c\u043d\u044b\u0439 \u043f\u043e \u043c\u043e\u0434\u0443\u043b\u044e \u044d\u043b\u0435\u043c\u0435\u043d\u0442 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0441\u0442\u043e\u043b\u0431\u0435\u0446 \u0438\u0437 \u043d\u0443\u043b\u0435\u0439.\",\n    2\n)\n\n# --- Change 2: remove the pre-existing \"_GoBack\" bookmark ---------------\n# Word keeps a single \"_GoBack\" bookmark that marks the location of the\n# last edit; it gets removed from its old spot and re-created at the new\n# edit location (see below).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Change 3: row 11 sentence rewording + relocate \"_GoBack\" -----------\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\n    \"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043d\u0435\u0433\u043e \u043a\u0430\u0436\u0434\u044b\u0439 \u0441\u0442\u043e\u043b\u0431\u0435\u0446, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0438\u0439 \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u043a\u0440\u0430\u0442\u043d\u044b\u0439 \u043f\u044f\u0442\u0438\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043c\u0430\u0441\u0441\u0438\u0432\u0430 \u043a\u0430\u0436\u0434\u044b\u0439 \u0441\u0442\u043e\u043b\u0431\u0435\u0446, \u0441\u043e\u0434\u0435\u0440\u0436\u0430\u0449\u0438\u0439 \u044d\u043b\u0435\u043c\u0435\u043d\u0442, \u043a\u0440\u0430\u0442\u043d\u044b\u0439 \u043f\u044f\u0442\u0438\",\n    2\n)\n\n# Find the freshly inserted word \"\u043c\u0430\u0441\u0441\u0438\u0432\u0430\" within this same sentence so we\n# can drop the \"_GoBack\" bookmark immediately after it.\n$range3 = $d.Content\n$found3 = $range3.Find.Execute(\"\u0423\u0434\u0430\u043b\u0438\u0442\u044c \u0438\u0437 \u043c\u0430\u0441\u0441\u0438\u0432\u0430\")\nif ($found3) {\n    $range3.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $range3)\n}\n"}
